# Updates the cryptos price/volume table with a fresh snapshot of values
# (and re-sorts a handful of coins whose ranking shifted), matching the
# "Updated cryptos list ... with GitHub Actions" commit.
#
# All of B/C/D/E are stored as text cells (not numbers), so for any D-column
# value that looks numeric (e.g. "1.007", "0.1046") we briefly force the
# cell to Text format before assigning it, otherwise Excel's COM layer will
# silently coerce the string into a float/int and mangle things like
# trailing zeros (1.170 -> 1.17) or tiny magnitudes (scientific notation).
# The number format is restored immediately afterwards so the cell's
# logical format stays "General", same as before the edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origFormat = $range.NumberFormat
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.NumberFormat = $origFormat
}

$ws.Range('D2').Value = '30.211.40'
$ws.Range('E2').Value = '  +1.08%  '
$ws.Range('D3').Value = '2.084.01'
$ws.Range('E3').Value = '  -1.60%  '
Set-TextValue $ws.Range('D4') '1.007'
$ws.Range('E4').Value = '  +0.10%  '
Set-TextValue $ws.Range('D5') '338.11'
$ws.Range('E5').Value = '  -2.85%  '
$ws.Range('E6').Value = '  -0.12%  '
Set-TextValue $ws.Range('D7') '0.5266'
$ws.Range('E7').Value = '  +1.47%  '
Set-TextValue $ws.Range('D8') '0.4367'
$ws.Range('E8').Value = '  -2.13%  '
Set-TextValue $ws.Range('D9') '54.84'
$ws.Range('E9').Value = '  +1.42%  '
Set-TextValue $ws.Range('D10') '0.09308'
$ws.Range('E10').Value = '  -0.64%  '
Set-TextValue $ws.Range('D11') '1.170'
$ws.Range('E11').Value = '  -0.91%  '
Set-TextValue $ws.Range('D12') '24.44'
$ws.Range('E12').Value = '  -3.03%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '2.104.89'
$ws.Range('E13').Value = '  -1.76%  '
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D14') '8.457'
$ws.Range('E14').Value = '  +0.55%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D15') '6.837'
$ws.Range('E15').Value = '  -0.24%  '
Set-TextValue $ws.Range('D16') '101.25'
$ws.Range('E16').Value = '  -1.17%  '
Set-TextValue $ws.Range('D18') '1.007'
$ws.Range('E18').Value = '  +0.03%  '
Set-TextValue $ws.Range('D19') '20.92'
$ws.Range('E19').Value = '  -2.96%  '
Set-TextValue $ws.Range('D20') '0.06716'
$ws.Range('E20').Value = '  +0.62%  '
Set-TextValue $ws.Range('D21') '6.289'
$ws.Range('E21').Value = '  -0.28%  '
$ws.Range('E22').Value = '  -0.26%  '
$ws.Range('D23').Value = '30.227.25'
$ws.Range('E23').Value = '  +1.00%  '
Set-TextValue $ws.Range('D24') '12.37'
$ws.Range('E24').Value = '  -2.72%  '
Set-TextValue $ws.Range('D25') '2.321'
$ws.Range('E25').Value = '  -0.24%  '
Set-TextValue $ws.Range('D26') '21.76'
$ws.Range('E26').Value = '  -1.60%  '
$ws.Range('B27').Value = 'Monero'
$ws.Range('C27').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws.Range('D27') '162.45'
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D28') '6.789'
$ws.Range('E28').Value = '  +4.34%  '
Set-TextValue $ws.Range('D29') '2.479'
$ws.Range('E29').Value = '  -3.59%  '
Set-TextValue $ws.Range('D30') '133.42'
$ws.Range('E30').Value = '  -0.43%  '
Set-TextValue $ws.Range('D31') '1.124'
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue $ws.Range('D32') '0.1046'
$ws.Range('E32').Value = '  -0.96%  '
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D33') '1.653'
$ws.Range('E33').Value = '  -8.32%  '
Set-TextValue $ws.Range('D34') '6.240'
$ws.Range('E34').Value = '  -0.24%  '
Set-TextValue $ws.Range('D35') '3.917'
$ws.Range('E35').Value = '  -1.44%  '
Set-TextValue $ws.Range('D36') '0.02611'
$ws.Range('E36').Value = '  +0.24%  '
Set-TextValue $ws.Range('D37') '9.820'
$ws.Range('E37').Value = '  -9.39%  '
Set-TextValue $ws.Range('D38') '0.06727'
$ws.Range('E38').Value = '  -1.37%  '
$ws.Range('B39').Value = 'TrustWalletToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D39') '1.343'
$ws.Range('E39').Value = '  +0.05%  '
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue $ws.Range('D40') '0.6927'
$ws.Range('E40').Value = '  -1.43%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue $ws.Range('D41') '12.50'
$ws.Range('E41').Value = '  -1.50%  '
Set-TextValue $ws.Range('D42') '0.2200'
$ws.Range('E42').Value = '  -2.04%  '
Set-TextValue $ws.Range('D43') '0.6713'
$ws.Range('E43').Value = '  -1.97%  '
Set-TextValue $ws.Range('D44') '2.366'
$ws.Range('E44').Value = '  +0.40%  '
Set-TextValue $ws.Range('D45') '14.31'
$ws.Range('E45').Value = '  -1.04%  '
Set-TextValue $ws.Range('D46') '1.005'
$ws.Range('E46').Value = '  -0.07%  '
Set-TextValue $ws.Range('D47') '1.293'
$ws.Range('E47').Value = '  +3.61%  '
Set-TextValue $ws.Range('D48') '3.629'
$ws.Range('E48').Value = '  -0.14%  '
Set-TextValue $ws.Range('D49') '0.00000000342'
$ws.Range('E49').Value = '  -4.80%  '
Set-TextValue $ws.Range('D50') '1.203'
$ws.Range('E50').Value = '  +1.93%  '
$ws.Range('E51').Value = '  -1.21%  '
